# Apply updated TPM-derived values to Cxcl12-Cd4 LR-pairs sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = [double]"146.620486"
$ws.Cells.Item(2, 8).Value = [double]"439.861458"
$ws.Cells.Item(2, 9).Value = [double]"0.3983053592962091"
$ws.Cells.Item(2, 10).Value = [double]"0.3983053592962091"
$ws.Cells.Item(2, 11).Value = [double]"3"
$ws.Cells.Item(2, 12).Value = [double]"1"
$ws.Cells.Item(2, 13).Value = [double]"2.441151666666667"
$ws.Cells.Item(2, 14).Value = [double]"7.323455"
$ws.Cells.Item(2, 15).Value = [double]"0.1131710230270566"
$ws.Cells.Item(2, 16).Value = [double]"0.1131710230270566"
$ws.Cells.Item(2, 17).Value = [double]"357.9228437663767"
$ws.Cells.Item(2, 18).Value = [double]"3221.30559389739"
$ws.Cells.Item(2, 19).Value = [double]"0.04507662498871132"
$ws.Cells.Item(2, 20).Value = [double]"0.04507662498871132"
$ws.Cells.Item(3, 7).Value = [double]"146.620486"
$ws.Cells.Item(3, 8).Value = [double]"439.861458"
$ws.Cells.Item(3, 9).Value = [double]"0.3983053592962091"
$ws.Cells.Item(3, 10).Value = [double]"0.3983053592962091"
$ws.Cells.Item(3, 15).Value = [double]"0.3310039188305578"
$ws.Cells.Item(3, 16).Value = [double]"0.3310039188305577"
$ws.Cells.Item(3, 17).Value = [double]"1046.856878702279"
$ws.Cells.Item(3, 18).Value = [double]"9421.711908320507"
$ws.Cells.Item(3, 19).Value = [double]"0.1318406348182586"
$ws.Cells.Item(3, 20).Value = [double]"0.1318406348182585"
$ws.Cells.Item(4, 7).Value = [double]"146.620486"
$ws.Cells.Item(4, 8).Value = [double]"439.861458"
$ws.Cells.Item(4, 9).Value = [double]"0.3983053592962091"
$ws.Cells.Item(4, 10).Value = [double]"0.3983053592962091"
$ws.Cells.Item(4, 11).Value = [double]"1"
$ws.Cells.Item(4, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(4, 13).Value = [double]"0.1123943333333333"
$ws.Cells.Item(4, 14).Value = [double]"0.337183"
$ws.Cells.Item(4, 15).Value = [double]"0.005210565922413945"
$ws.Cells.Item(4, 16).Value = [double]"0.005210565922413945"
$ws.Cells.Item(4, 17).Value = [double]"16.47931177697933"
$ws.Cells.Item(4, 18).Value = [double]"148.313805992814"
$ws.Cells.Item(4, 19).Value = [double]"0.00207539633186367"
$ws.Cells.Item(4, 20).Value = [double]"0.00207539633186367"
$ws.Cells.Item(5, 7).Value = [double]"146.620486"
$ws.Cells.Item(5, 8).Value = [double]"439.861458"
$ws.Cells.Item(5, 9).Value = [double]"0.3983053592962091"
$ws.Cells.Item(5, 10).Value = [double]"0.3983053592962091"
$ws.Cells.Item(5, 13).Value = [double]"11.877011"
$ws.Cells.Item(5, 14).Value = [double]"35.631033"
$ws.Cells.Item(5, 15).Value = [double]"0.5506144922199717"
$ws.Cells.Item(5, 16).Value = [double]"0.5506144922199717"
$ws.Cells.Item(5, 17).Value = [double]"1741.413125047346"
$ws.Cells.Item(5, 18).Value = [double]"15672.71812542611"
$ws.Cells.Item(5, 19).Value = [double]"0.2193127031573756"
$ws.Cells.Item(5, 20).Value = [double]"0.2193127031573756"
$ws.Cells.Item(6, 9).Value = [double]"0.534552907532962"
$ws.Cells.Item(6, 10).Value = [double]"0.5345529075329621"
$ws.Cells.Item(6, 11).Value = [double]"3"
$ws.Cells.Item(6, 12).Value = [double]"1"
$ws.Cells.Item(6, 13).Value = [double]"2.441151666666667"
$ws.Cells.Item(6, 14).Value = [double]"7.323455"
$ws.Cells.Item(6, 15).Value = [double]"0.1131710230270566"
$ws.Cells.Item(6, 16).Value = [double]"0.1131710230270566"
$ws.Cells.Item(6, 17).Value = [double]"480.3568225791728"
$ws.Cells.Item(6, 18).Value = [double]"4323.211403212556"
$ws.Cells.Item(6, 19).Value = [double]"0.0604958994075929"
$ws.Cells.Item(6, 20).Value = [double]"0.0604958994075929"
$ws.Cells.Item(7, 9).Value = [double]"0.534552907532962"
$ws.Cells.Item(7, 10).Value = [double]"0.5345529075329621"
$ws.Cells.Item(7, 15).Value = [double]"0.3310039188305578"
$ws.Cells.Item(7, 16).Value = [double]"0.3310039188305577"
$ws.Cells.Item(7, 19).Value = [double]"0.1769391072156792"
$ws.Cells.Item(7, 20).Value = [double]"0.1769391072156792"
$ws.Cells.Item(8, 9).Value = [double]"0.534552907532962"
$ws.Cells.Item(8, 10).Value = [double]"0.5345529075329621"
$ws.Cells.Item(8, 11).Value = [double]"1"
$ws.Cells.Item(8, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(8, 13).Value = [double]"0.1123943333333333"
$ws.Cells.Item(8, 14).Value = [double]"0.337183"
$ws.Cells.Item(8, 15).Value = [double]"0.005210565922413945"
$ws.Cells.Item(8, 16).Value = [double]"0.005210565922413945"
$ws.Cells.Item(8, 17).Value = [double]"22.11635826364922"
$ws.Cells.Item(8, 18).Value = [double]"199.047224372843"
$ws.Cells.Item(8, 19).Value = [double]"0.002785323163718545"
$ws.Cells.Item(8, 20).Value = [double]"0.002785323163718545"
$ws.Cells.Item(9, 9).Value = [double]"0.534552907532962"
$ws.Cells.Item(9, 10).Value = [double]"0.5345529075329621"
$ws.Cells.Item(9, 13).Value = [double]"11.877011"
$ws.Cells.Item(9, 14).Value = [double]"35.631033"
$ws.Cells.Item(9, 15).Value = [double]"0.5506144922199717"
$ws.Cells.Item(9, 16).Value = [double]"0.5506144922199717"
$ws.Cells.Item(9, 17).Value = [double]"2337.09496366041"
$ws.Cells.Item(9, 18).Value = [double]"21033.8546729437"
$ws.Cells.Item(9, 19).Value = [double]"0.2943325777459714"
$ws.Cells.Item(9, 20).Value = [double]"0.2943325777459714"
$ws.Cells.Item(10, 7).Value = [double]"24.174389"
$ws.Cells.Item(10, 8).Value = [double]"72.523167"
$ws.Cells.Item(10, 9).Value = [double]"0.0656715098899026"
$ws.Cells.Item(10, 10).Value = [double]"0.0656715098899026"
$ws.Cells.Item(10, 11).Value = [double]"3"
$ws.Cells.Item(10, 12).Value = [double]"1"
$ws.Cells.Item(10, 13).Value = [double]"2.441151666666667"
$ws.Cells.Item(10, 14).Value = [double]"7.323455"
$ws.Cells.Item(10, 15).Value = [double]"0.1131710230270566"
$ws.Cells.Item(10, 16).Value = [double]"0.1131710230270566"
$ws.Cells.Item(10, 17).Value = [double]"59.01334999799834"
$ws.Cells.Item(10, 18).Value = [double]"531.120149981985"
$ws.Cells.Item(10, 19).Value = [double]"0.007432111957971742"
$ws.Cells.Item(10, 20).Value = [double]"0.007432111957971741"
$ws.Cells.Item(11, 7).Value = [double]"24.174389"
$ws.Cells.Item(11, 8).Value = [double]"72.523167"
$ws.Cells.Item(11, 9).Value = [double]"0.0656715098899026"
$ws.Cells.Item(11, 10).Value = [double]"0.0656715098899026"
$ws.Cells.Item(11, 15).Value = [double]"0.3310039188305578"
$ws.Cells.Item(11, 16).Value = [double]"0.3310039188305577"
$ws.Cells.Item(11, 17).Value = [double]"172.6029295324713"
$ws.Cells.Item(11, 18).Value = [double]"1553.426365792242"
$ws.Cells.Item(11, 19).Value = [double]"0.02173752712907749"
$ws.Cells.Item(11, 20).Value = [double]"0.02173752712907749"
$ws.Cells.Item(12, 7).Value = [double]"24.174389"
$ws.Cells.Item(12, 8).Value = [double]"72.523167"
$ws.Cells.Item(12, 9).Value = [double]"0.0656715098899026"
$ws.Cells.Item(12, 10).Value = [double]"0.0656715098899026"
$ws.Cells.Item(12, 11).Value = [double]"1"
$ws.Cells.Item(12, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(12, 13).Value = [double]"0.1123943333333333"
$ws.Cells.Item(12, 14).Value = [double]"0.337183"
$ws.Cells.Item(12, 15).Value = [double]"0.005210565922413945"
$ws.Cells.Item(12, 16).Value = [double]"0.005210565922413945"
$ws.Cells.Item(12, 17).Value = [double]"2.717064335395667"
$ws.Cells.Item(12, 18).Value = [double]"24.453579018561"
$ws.Cells.Item(12, 19).Value = [double]"0.0003421857315057969"
$ws.Cells.Item(12, 20).Value = [double]"0.0003421857315057969"
$ws.Cells.Item(13, 7).Value = [double]"24.174389"
$ws.Cells.Item(13, 8).Value = [double]"72.523167"
$ws.Cells.Item(13, 9).Value = [double]"0.0656715098899026"
$ws.Cells.Item(13, 10).Value = [double]"0.0656715098899026"
$ws.Cells.Item(13, 13).Value = [double]"11.877011"
$ws.Cells.Item(13, 14).Value = [double]"35.631033"
$ws.Cells.Item(13, 15).Value = [double]"0.5506144922199717"
$ws.Cells.Item(13, 16).Value = [double]"0.5506144922199717"
$ws.Cells.Item(13, 17).Value = [double]"287.119484071279"
$ws.Cells.Item(13, 18).Value = [double]"2584.075356641511"
$ws.Cells.Item(13, 19).Value = [double]"0.03615968507134757"
$ws.Cells.Item(13, 20).Value = [double]"0.03615968507134757"
$ws.Cells.Item(14, 5).Value = [double]"3"
$ws.Cells.Item(14, 6).Value = [double]"1"
$ws.Cells.Item(14, 7).Value = [double]"0.541205"
$ws.Cells.Item(14, 8).Value = [double]"1.623615"
$ws.Cells.Item(14, 9).Value = [double]"0.001470223280926138"
$ws.Cells.Item(14, 10).Value = [double]"0.001470223280926138"
$ws.Cells.Item(14, 11).Value = [double]"3"
$ws.Cells.Item(14, 12).Value = [double]"1"
$ws.Cells.Item(14, 13).Value = [double]"2.441151666666667"
$ws.Cells.Item(14, 14).Value = [double]"7.323455"
$ws.Cells.Item(14, 15).Value = [double]"0.1131710230270566"
$ws.Cells.Item(14, 16).Value = [double]"0.1131710230270566"
$ws.Cells.Item(14, 17).Value = [double]"1.321163487758334"
$ws.Cells.Item(14, 18).Value = [double]"11.890471389825"
$ws.Cells.Item(14, 19).Value = [double]"0.0001663866727806066"
$ws.Cells.Item(14, 20).Value = [double]"0.0001663866727806066"
$ws.Cells.Item(15, 5).Value = [double]"3"
$ws.Cells.Item(15, 6).Value = [double]"1"
$ws.Cells.Item(15, 7).Value = [double]"0.541205"
$ws.Cells.Item(15, 8).Value = [double]"1.623615"
$ws.Cells.Item(15, 9).Value = [double]"0.001470223280926138"
$ws.Cells.Item(15, 10).Value = [double]"0.001470223280926138"
$ws.Cells.Item(15, 15).Value = [double]"0.3310039188305578"
$ws.Cells.Item(15, 16).Value = [double]"0.3310039188305577"
$ws.Cells.Item(15, 17).Value = [double]"3.864154269943334"
$ws.Cells.Item(15, 18).Value = [double]"34.77738842949"
$ws.Cells.Item(15, 19).Value = [double]"0.0004866496675424717"
$ws.Cells.Item(15, 20).Value = [double]"0.0004866496675424716"
$ws.Cells.Item(16, 5).Value = [double]"3"
$ws.Cells.Item(16, 6).Value = [double]"1"
$ws.Cells.Item(16, 7).Value = [double]"0.541205"
$ws.Cells.Item(16, 8).Value = [double]"1.623615"
$ws.Cells.Item(16, 9).Value = [double]"0.001470223280926138"
$ws.Cells.Item(16, 10).Value = [double]"0.001470223280926138"
$ws.Cells.Item(16, 11).Value = [double]"1"
$ws.Cells.Item(16, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(16, 13).Value = [double]"0.1123943333333333"
$ws.Cells.Item(16, 14).Value = [double]"0.337183"
$ws.Cells.Item(16, 15).Value = [double]"0.005210565922413945"
$ws.Cells.Item(16, 16).Value = [double]"0.005210565922413945"
$ws.Cells.Item(16, 17).Value = [double]"0.06082837517166667"
$ws.Cells.Item(16, 18).Value = [double]"0.547455376545"
$ws.Cells.Item(16, 19).Value = [double]"7.660695325933359E-06"
$ws.Cells.Item(16, 20).Value = [double]"7.660695325933359E-06"
$ws.Cells.Item(17, 5).Value = [double]"3"
$ws.Cells.Item(17, 6).Value = [double]"1"
$ws.Cells.Item(17, 7).Value = [double]"0.541205"
$ws.Cells.Item(17, 8).Value = [double]"1.623615"
$ws.Cells.Item(17, 9).Value = [double]"0.001470223280926138"
$ws.Cells.Item(17, 10).Value = [double]"0.001470223280926138"
$ws.Cells.Item(17, 13).Value = [double]"11.877011"
$ws.Cells.Item(17, 14).Value = [double]"35.631033"
$ws.Cells.Item(17, 15).Value = [double]"0.5506144922199717"
$ws.Cells.Item(17, 16).Value = [double]"0.5506144922199717"
$ws.Cells.Item(17, 17).Value = [double]"6.427897738255001"
$ws.Cells.Item(17, 18).Value = [double]"57.85107964429501"
$ws.Cells.Item(17, 19).Value = [double]"0.0008095262452771263"
$ws.Cells.Item(17, 20).Value = [double]"0.0008095262452771263"
